# Apply the layup-data edits described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("layer 1"): material_id 2 -> 1, keep orientation 0, clear thickness value
$ws.Range("B2").Value = 1

# Row 3 ("layer 2"): material_id 2 -> 1, orientation 0 -> 45, clear thickness value
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 45

# The "thickness" column data (D1 header + D2/D3 values) is removed; the
# empty D1 header cell disappears entirely, while D2/D3 keep their number
# formatting but become blank.
$ws.Range("D1").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("D3").ClearContents()

# Move the active selection/cursor to C3.
$ws.Range("C3").Select()
